# Fix bugs in templates: rename "AVISO" header to "INCIDENTE" on the
# EDIFICIOS sheet (cell G1), matching the style already used by the
# adjacent header cells (D1/E1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EDIFICIOS")

$g1 = $ws.Range("G1")
$g1.Value = "INCIDENTE"
$g1.Font.Name = "Arial"
$g1.Font.Color = 16777215
$g1.Font.Bold = $true
